# New crime data collected - update CompStat_1 weekly report
# (Volume/Number header, reporting week dates, and the crime-complaint
# statistics table for rows 15-30.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header: "Volume 30   Number  24" -> "Volume 30   Number  25" ---
$ws.Range("A8").Value = "Volume 30   Number  25"

# --- Header: reporting week "6/12/2023 .. 6/18/2023" -> "6/19/2023 .. 6/25/2023" ---
$ws.Range("C9").Value = "Report Covering the Week  6/19/2023  Through  6/25/2023"

# --- Crime complaints table (rows 15-30) ---

# Row 15
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 2
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 8
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -20
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 700
$ws.Range("N15").Value = 0

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 120
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = -9.090909090909
$ws.Range("L16").Value = 36.363636363636
$ws.Range("M16").Value = 71.428571428571
$ws.Range("N16").Value = -85.294117647058

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -53.846153846153
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = -1.639344262295
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 122.222222222222
$ws.Range("N17").Value = -24.050632911392

# Row 18
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -55
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = -32.089552238806
$ws.Range("L18").Value = 62.5
$ws.Range("M18").Value = 5.813953488372
$ws.Range("N18").Value = -77.804878048780

# Row 19
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = 18.181818181818
$ws.Range("F19").Value = 92
$ws.Range("G19").Value = 102
$ws.Range("H19").Value = -9.803921568627
$ws.Range("I19").Value = 551
$ws.Range("J19").Value = 584
$ws.Range("K19").Value = -5.650684931506
$ws.Range("L19").Value = 81.25
$ws.Range("M19").Value = 9.325396825396
$ws.Range("N19").Value = -68.442153493699

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 35
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = 12.903225806451
$ws.Range("L20").Value = 75
$ws.Range("M20").Value = 133.333333333333
$ws.Range("N20").Value = -91.879350348027

# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 13.157894736842
$ws.Range("F21").Value = 131
$ws.Range("G21").Value = 150
$ws.Range("H21").Value = -12.666666666666
$ws.Range("I21").Value = 805
$ws.Range("J21").Value = 886
$ws.Range("K21").Value = -9.142212189616
$ws.Range("L21").Value = 62.955465587044
$ws.Range("M21").Value = 20.508982035928
$ws.Range("N21").Value = -73.948220064724

# Row 22
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -66.666666666666
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = -37.5
$ws.Range("I22").Value = 35
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = -35.185185185185
$ws.Range("L22").Value = 9.375
$ws.Range("M22").Value = 2.941176470588

# Row 24
$ws.Range("C24").Value = 66
$ws.Range("D24").Value = 85
$ws.Range("E24").Value = -22.352941176470
$ws.Range("F24").Value = 299
$ws.Range("H24").Value = -12.316715542522
$ws.Range("I24").Value = 1829
$ws.Range("J24").Value = 1894
$ws.Range("K24").Value = -3.431890179514
$ws.Range("L24").Value = 91.317991631799
$ws.Range("M24").Value = 121.428571428571

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 173
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 8.125
$ws.Range("L25").Value = 17.687074829932
$ws.Range("M25").Value = 41.803278688524

# Row 26
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 10
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = 11.111111111111

# Row 27
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 42
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = -16
$ws.Range("L27").Value = 13.513513513513

# Row 30
$ws.Range("L30").Value = -42.857142857142
